$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.192.55"
$ws.Range("E2").Value = "  -3.65%  "

$ws.Range("D3").Value = "3.483.36"
$ws.Range("E3").Value = "  -5.40%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "607.12"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -6.87%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.94"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -8.22%  "

$ws.Range("D7").Value = "3.481.97"
$ws.Range("E7").Value = "  -5.37%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.478"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -3.96%  "

$ws.Range("E10").Value = "  -5.41%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.91"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -3.46%  "

$ws.Range("E12").Value = "  -4.66%  "

$ws.Range("E13").Value = "  -6.85%  "

$ws.Range("D14").Value = "4.073.02"
$ws.Range("E14").Value = "  -5.35%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "31.34"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -4.15%  "

$ws.Range("D16").Value = "3.489.07"
$ws.Range("E16").Value = "  -4.94%  "

$ws.Range("D17").Value = "67.106.47"
$ws.Range("E17").Value = "  -3.78%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.116"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.91%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.41"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.70%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.96"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -6.07%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "445.89"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -5.23%  "

$ws.Range("E22").Value = "  -13.07%  "

$ws.Range("E23").Value = "  -5.21%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "76.99"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -3.56%  "

$ws.Range("E25").Value = "  +0.10%  "

$ws.Range("E26").Value = "  -0.28%  "

$ws.Range("D27").Value = "3.622.47"
$ws.Range("E27").Value = "  -5.34%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.09"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -9.51%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.26"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -5.63%  "

$ws.Range("E30").Value = "  -4.88%  "

$ws.Range("E31").Value = "  -7.71%  "

$ws.Range("E32").Value = "  +0.38%  "

$ws.Range("E33").Value = "  -1.74%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.69"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -3.94%  "

$ws.Range("E35").Value = "  -5.34%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.84"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -7.32%  "

$ws.Range("D37").Value = "3.477.49"
$ws.Range("E37").Value = "  -5.40%  "

$ws.Range("E38").Value = "  -5.06%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.08%  "

$ws.Range("E41").Value = "  +0.70%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "170.37"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -4.40%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0869"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.99%  "

$ws.Range("E44").Value = "  -7.94%  "

$ws.Range("E45").Value = "  -5.29%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "45.35"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.15%  "

$ws.Range("E47").Value = "  +0.96%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "26.54"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -8.77%  "

$ws.Range("E49").Value = "  -9.24%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.51"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -4.42%  "

$ws.Range("E51").Value = "  -4.04%  "
